$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire column D ("下一次充电开始时间"), shifting columns E:F left to D:E
$ws.Range("D1").EntireColumn.Delete()

# The AutoFilter / _FilterDatabase defined name still references the old
# F column; update it to reflect the new right-most column (E) after the
# column deletion.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "sheet1!_FilterDatabase") {
        $n.RefersTo = "=sheet1!`$A`$1:`$E`$73"
    }
}

# Update the active selection to match the post-edit workbook state
$ws.Range("G7").Select()
